$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jobDescription = "Please find the Job Description (JD) below for your reference.`nJob Role: Landing Gear System Engineer`nJob Description:`n" +
"• Been able to analyse landing gear systems, ensuring they meet performance, safety, and regulatory requirements.`n" +
"•Collaborate with interdisciplinary teams to seamlessly integrate landing gear mechanisms into the overall aircraft design.`n" +
"•Design of hydraulic or electric actuation systems for smooth and reliable landing gear extension and retraction.`n" +
"•Ensure landing gear systems adhere to aviation safety standards and regulations, conducting thorough risk assessments.`n" +
"•Conduct rigorous testing protocols to validate the functionality and reliability of landing gear systems under various conditions.`n" +
"•Continuously strive to optimize landing gear performance in terms of weight, efficiency, and overall functionality.`n" +
"•Troubleshoot and address any issues related to landing gear operation, collaborating with other engineering disciplines as needed.`n" +
"•Maintain comprehensive documentation of design specifications, test results, and compliance records.`n" +
"•Work closely with aerodynamics, structural, and systems engineers to ensure a holistic approach to aircraft design and performance.`n" +
"•Participate in the development and testing of models to validate design concepts and identify areas for improvement.`n" +
"•Effectively communicate design concepts, progress, and challenges to team members, project managers, and other stakeholders.`n" +
"Qualifications/Required Skills:`n" +
"• Bachelor’s or Master's degree in Mechanical Engineering or any other relavent field`n" +
"• 7 yr - 15 years of professional experience in Aircraft Landing Gear systems.`n" +
"• Basic awareness of Aircraft System Qualification & Certification.`n" +
"•Process optimisation and organisational skills.`n" +
"•Planning & reporting experience.`n" +
"•Project management experience is desired but not necessary.`n" +
"•Creating & analysing test cases from lab and flight test.`n" +
"•Ability to support delivery of large complex projects on time.`n" +
"•Experience in writing technical specifications both at equipment & systems level.`n" +
"•Systems engineering experience; particularly involving Requirements Management (DOORS), Integration, Verification & Validation.`n" +
"•Aerospace (Landing Gear Systems) environment background.`n" +
"•Proficiency in Hydraulic & Electromechanical systems for landing gear deployment. `n" +
"•Experience in the use of one or more of the following would be advantageous: Catia V5, MSC ADAMS, MATLAB, Simulink or IBM DOORS."

$ws.Range("A3").Value = $jobDescription
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 7
